# Update "want-to-go" attendee counts (column F) on the "展览" and "全部类型"
# sheets to reflect the latest scrape, per the gh-pages output regeneration.

$wb = $excel.ActiveWorkbook

# Row (by sheet row number) -> new F value
$updates = @{
    2  = 379
    3  = 10877
    5  = 984
    6  = 196
    7  = 1346
    8  = 8335
    9  = 45
    11 = 616
    12 = 223
    14 = 3331
    16 = 330
    18 = 825
    19 = 134
    20 = 1077
    22 = 129
    23 = 1845
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
